$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "55÷8=6, 7" "61÷4=15, 1"
Replace-Text "73÷3=24, 1" "99÷8=12, 3"
Replace-Text "48÷5=9, 3" "12÷8=1, 4"
Replace-Text "66÷8=8, 2" "29÷4=7, 1"
Replace-Text "95÷9=10, 5" "48÷2=24, 0"
Replace-Text "15÷8=1, 7" "27÷8=3, 3"
Replace-Text "78÷3=26, 0" "14÷5=2, 4"
Replace-Text "16÷5=3, 1" "17÷4=4, 1"
Replace-Text "63÷4=15, 3" "26÷3=8, 2"
Replace-Text "64÷4=16, 0" "54÷7=7, 5"
Replace-Text "67÷5=13, 2" "20÷2=10, 0"
Replace-Text "21÷6=3, 3" "92÷4=23, 0"
Replace-Text "90÷7=12, 6" "51÷5=10, 1"
Replace-Text "45÷3=15, 0" "74÷8=9, 2"
Replace-Text "62÷5=12, 2" "46÷3=15, 1"
Replace-Text "53÷3=17, 2" "58÷7=8, 2"
Replace-Text "66÷6=11, 0" "65÷3=21, 2"
Replace-Text "57÷8=7, 1" "87÷8=10, 7"
Replace-Text "67÷7=9, 4" "66÷9=7, 3"
Replace-Text "46÷4=11, 2" "49÷5=9, 4"
Replace-Text "93÷4=23, 1" "77÷2=38, 1"
Replace-Text "35÷4=8, 3" "37÷8=4, 5"
Replace-Text "24÷8=3, 0" "66÷9=7, 3"
Replace-Text "90÷5=18, 0" "97÷7=13, 6"
Replace-Text "15÷9=1, 6" "50÷4=12, 2"
